$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Append a new row for issue #27 (keep the Issue ID as text, matching the
# rest of column A, which stores issue numbers as text even though they
# look numeric). Temporarily format A23 as Text so Excel doesn't coerce
# "27" into a number, then clear the format again so the cell keeps the
# sheet's default (unstyled) look, matching the other rows.
$row = 23
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "27"
$ws.Range("A$row").ClearFormats()
$ws.Range("B$row").Value = "[FEATURE REQUEST] <title>"
$ws.Range("C$row").Value = "open"
$ws.Range("D$row").Value = "2025-03-26T06:59:28Z"
$ws.Range("E$row").Value = "enhancement"
